$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("D2").Value = "亭洪路45号 百益上河城"
$ws1.Range("F2").Value = 5460
$ws1.Range("F3").Value = 595
$ws1.Range("F4").Value = 11932
$ws1.Range("F5").Value = 294
$ws1.Range("F7").Value = 176
$ws1.Range("F8").Value = 303
$ws1.Range("F9").Value = 1074

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("D4").Value = "亭洪路45号 百益上河城"
$ws4.Range("F4").Value = 5460
$ws4.Range("F5").Value = 595
$ws4.Range("F7").Value = 11932
$ws4.Range("F8").Value = 294
$ws4.Range("F10").Value = 176
$ws4.Range("F13").Value = 303
$ws4.Range("F14").Value = 1074
